# Applies the text corrections described by the commit diff.
# Uses TextRange.Characters(1, Length).Text = "..." rather than
# TextRange.Text = "..." directly: both update the run text, but going
# through Characters() on the full span keeps the paragraph's run intact
# as a single run and does not drop the shape's trailing empty paragraph,
# matching PowerPoint's normal in-place text-replace behavior.

$p = $ppt.ActivePresentation

function Set-ShapeText {
    param(
        [int]$SlideIndex,
        [int]$ShapeIndex,
        [string]$NewText
    )
    $tr = $p.Slides.Item($SlideIndex).Shapes.Item($ShapeIndex).TextFrame.TextRange
    $full = $tr.Characters(1, $tr.Length)
    $full.Text = $NewText
}

# 第三場 精神科にて -> 第三場 心療内科クリニックにて
Set-ShapeText 119 2 "第三場 心療内科クリニックにて"

# 口答期 -> 口唇期
Set-ShapeText 134 2 "声: そんな乱れた関係は あなたが未熟で まだ口唇期を乗り越えていない証拠です"

# MARÍA typo: "es había ido" -> "se había ido"
Set-ShapeText 149 1 "MARÍA. - Me di cuenta cuando se había ido, me lo quitó del joyero... "

# duplicated "に" typo: "ものにに泣いている" -> "ものに泣いている"
Set-ShapeText 51 2 "警部: 祖国の危機の最中 あなたは平手打のようなものに泣いている!"

# INSPECTOR typo: "perdodo" -> "perdono"
Set-ShapeText 55 1 "INSPECTOR, - Ande, ¡váyase! ¡Váyase de una vez y por ésta se lo perdono!"

# 第二場 法律相談所にて -> 第二場 弁護士事務所
Set-ShapeText 58 2 "第二場 弁護士事務所"

# 警部: なぜ? -> 警部: 何のために?
Set-ShapeText 6 2 "警部: 何のために?"

# 目撃者が明白な誓いのもと -> 目撃者の明白な誓いのもと
Set-ShapeText 76 2 "弁護士: それも 居合わせた目撃者の明白な誓いのもと!"

# 完全なる交が必要です -> 完全なる性交が必要です
Set-ShapeText 78 2 "弁護士: 姦通を事実と認めるには 完全なる性交が必要です"

# MARGARITA typo: "apar- tamento" -> "apartamento"
Set-ShapeText 87 1 "MARGARITA. - También tiene alquilado un apartamento donde van juntos al terminar el trabajo"

# 中人は家を出て行って -> 主人は家を出て行って
Set-ShapeText 95 2 "マルガリータ: それに 主人は家を出て行って 彼女と一緒にいるのですよ!"
